# use column indeces from sources.csv instead of hard coded
# The hard-coded "Category" column (E) is no longer produced by the
# source data, so remove it entirely and let the remaining training
# columns (F:M) shift left into E:L.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E:E").Delete()
